$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Range("D2") "286.90"
Set-TextValue $ws.Range("E2") "1.91%"
Set-TextValue $ws.Range("D3") "28.71"
Set-TextValue $ws.Range("E3") "3.96%"
Set-TextValue $ws.Range("D4") "5.075"
Set-TextValue $ws.Range("E4") "4.82%"
Set-TextValue $ws.Range("E5") "2.28%"
Set-TextValue $ws.Range("D6") "7.390"
Set-TextValue $ws.Range("E6") "4.03%"
Set-TextValue $ws.Range("D7") "3.404"
Set-TextValue $ws.Range("E7") "2.67%"
Set-TextValue $ws.Range("D8") "1.371"
Set-TextValue $ws.Range("E8") "7.81%"
Set-TextValue $ws.Range("D9") "0.9361"
Set-TextValue $ws.Range("E9") "2.03%"
Set-TextValue $ws.Range("D10") "0.1575"
Set-TextValue $ws.Range("E10") "1.48%"
Set-TextValue $ws.Range("D11") "0.06415"
Set-TextValue $ws.Range("E11") "-1.02%"
Set-TextValue $ws.Range("D12") "0.07627"
Set-TextValue $ws.Range("E12") "0.93%"
Set-TextValue $ws.Range("D13") "0.02935"
Set-TextValue $ws.Range("E13") "0.52%"
Set-TextValue $ws.Range("D14") "0.08979"
Set-TextValue $ws.Range("E14") "-0.45%"
Set-TextValue $ws.Range("D15") "0.001583"
Set-TextValue $ws.Range("E15") "0.43%"
Set-TextValue $ws.Range("D16") "0.04503"
Set-TextValue $ws.Range("E16") "2.18%"
Set-TextValue $ws.Range("D17") "0.0006450"
Set-TextValue $ws.Range("E17") "0.67%"
Set-TextValue $ws.Range("D18") "0.006266"
Set-TextValue $ws.Range("E18") "4.45%"
Set-TextValue $ws.Range("D19") "3.446"
Set-TextValue $ws.Range("E19") "-1.62%"
Set-TextValue $ws.Range("D20") "2.249"
Set-TextValue $ws.Range("E20") "1.05%"
Set-TextValue $ws.Range("D21") "0.3215"
Set-TextValue $ws.Range("E22") "-3.28%"
Set-TextValue $ws.Range("D23") "4.105"
Set-TextValue $ws.Range("D24") "0.1552"
Set-TextValue $ws.Range("E24") "4.41%"
Set-TextValue $ws.Range("D25") "0.001189"
Set-TextValue $ws.Range("E25") "1.93%"
Set-TextValue $ws.Range("D26") "0.004138"
Set-TextValue $ws.Range("E26") "-5.82%"
Set-TextValue $ws.Range("E27") "6.33%"
Set-TextValue $ws.Range("D28") "0.0001617"
Set-TextValue $ws.Range("E28") "-1.39%"
Set-TextValue $ws.Range("D40") "0.04215"
Set-TextValue $ws.Range("E40") "2.59%"
Set-TextValue $ws.Range("D41") "0.006750"
Set-TextValue $ws.Range("E41") "7.34%"
Set-TextValue $ws.Range("E42") "-10.67%"
Set-TextValue $ws.Range("D43") "0.001980"
Set-TextValue $ws.Range("E43") "-4.91%"
Set-TextValue $ws.Range("E44") "4.63%"
Set-TextValue $ws.Range("D45") "0.00005567"
Set-TextValue $ws.Range("E45") "1.46%"
Set-TextValue $ws.Range("D46") "1.970"
Set-TextValue $ws.Range("E46") "21.02%"
Set-TextValue $ws.Range("D47") "0.01306"
Set-TextValue $ws.Range("E47") "-28.60%"
